$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 314.07144
$ws.Range("I96").Value = 309.44446
$ws.Range("J96").Value = 322.4
$ws.Range("K96").Value = 928.33338
$ws.Range("L96").Value = 967.1999999999999
$ws.Range("M96").Value = 444.66662
$ws.Range("N96").Value = -3713.2

$ws.Range("H99").Value = 388
$ws.Range("J99").Value = 600
$ws.Range("L99").Value = 1800
$ws.Range("N99").Value = -4796

$ws.Range("H100").Value = 35716410
$ws.Range("I100").Value = 1433.3334
$ws.Range("J100").Value = 45456860
$ws.Range("K100").Value = 1433.3334
$ws.Range("L100").Value = 45456860
$ws.Range("M100").Value = -892.3334
$ws.Range("N100").Value = -45457942

$ws.Range("H112").Value = 2595.6667
$ws.Range("I112").Value = 350
$ws.Range("J112").Value = 3044.8
$ws.Range("K112").Value = 1050
$ws.Range("L112").Value = 9134.400000000001
$ws.Range("M112").Value = 58
$ws.Range("N112").Value = -11350.4

$ws.Range("H116").Value = 5200.6562
$ws.Range("I116").Value = 5056.4116
$ws.Range("J116").Value = 5364.1333
$ws.Range("K116").Value = 5056.4116
$ws.Range("L116").Value = 5364.1333
$ws.Range("M116").Value = -1614.4116
$ws.Range("N116").Value = -12248.1333

$ws.Range("H132").Value = 3541.2812
$ws.Range("I132").Value = 3518.3044
$ws.Range("J132").Value = 3600
$ws.Range("K132").Value = 10554.9132
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -8024.913199999999
$ws.Range("N132").Value = -15860

$ws.Range("H137").Value = 2831301
$ws.Range("I137").Value = 1316814.5
$ws.Range("K137").Value = 3950443.5
$ws.Range("M137").Value = -3947893.5

$ws.Range("H138").Value = 313491.66
$ws.Range("I138").Value = 4297.3335
$ws.Range("J138").Value = 363631.28
$ws.Range("K138").Value = 12892.0005
$ws.Range("L138").Value = 1090893.84
$ws.Range("M138").Value = -7752.000499999998
$ws.Range("N138").Value = -1101173.84

$ws.Range("H139").Value = 51375
$ws.Range("J139").Value = 51375
$ws.Range("L139").Value = 51375
$ws.Range("N139").Value = -61655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2621581
$ws.Range("I32").Value = 4364.46
$ws.Range("K32").Value = 4364.46
$ws.Range("M32").Value = -4077.46

$ws.Range("H61").Value = 1519.0667
$ws.Range("I61").Value = 1570.7858
$ws.Range("K61").Value = 1570.7858
$ws.Range("M61").Value = -1358.7858

$ws.Range("H74").Value = 967.3555
$ws.Range("I74").Value = 973.5135
$ws.Range("J74").Value = 938.875
$ws.Range("K74").Value = 973.5135
$ws.Range("L74").Value = 938.875
$ws.Range("M74").Value = -99.51350000000002
$ws.Range("N74").Value = -2686.875

$ws.Range("H77").Value = 967.3555
$ws.Range("I77").Value = 973.5135
$ws.Range("J77").Value = 938.875
$ws.Range("K77").Value = 4867.5675
$ws.Range("L77").Value = 4694.375
$ws.Range("M77").Value = -499.5675000000001
$ws.Range("N77").Value = -13430.375

$ws.Range("H132").Value = 77772.45
$ws.Range("I132").Value = 93879.22
$ws.Range("K132").Value = 281637.66
$ws.Range("M132").Value = -279107.66

$ws.Range("H134").Value = 32880.777
$ws.Range("J134").Value = 32880.777
$ws.Range("L134").Value = 32880.777
$ws.Range("N134").Value = -43020.777

$ws.Range("H136").Value = 1519.0667
$ws.Range("I136").Value = 1570.7858
$ws.Range("K136").Value = 4712.357400000001
$ws.Range("M136").Value = -2162.357400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2400.4167
$ws.Range("I105").Value = 1529.8334
$ws.Range("J105").Value = 4141.5835
$ws.Range("K105").Value = 1529.8334
$ws.Range("L105").Value = 4141.5835
$ws.Range("M105").Value = 217.1666
$ws.Range("N105").Value = -7635.5835

$ws.Range("H107").Value = 6705.625
$ws.Range("I107").Value = 7456.4165
$ws.Range("J107").Value = 4453.25
$ws.Range("K107").Value = 7456.4165
$ws.Range("L107").Value = 4453.25
$ws.Range("M107").Value = -5536.4165
$ws.Range("N107").Value = -8293.25

$ws.Range("H134").Value = 48021.516
$ws.Range("I134").Value = 55488.164
$ws.Range("J134").Value = 2392
$ws.Range("K134").Value = 166464.492
$ws.Range("L134").Value = 7176
$ws.Range("M134").Value = -163929.492
$ws.Range("N134").Value = -12246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1682.9454
$ws.Range("I31").Value = 1278
$ws.Range("J31").Value = 2588.1177
$ws.Range("K31").Value = 1278
$ws.Range("L31").Value = 2588.1177
$ws.Range("M31").Value = -983
$ws.Range("N31").Value = -3178.1177

$ws.Range("H34").Value = 1682.9454
$ws.Range("I34").Value = 1278
$ws.Range("J34").Value = 2588.1177
$ws.Range("K34").Value = 1278
$ws.Range("L34").Value = 2588.1177
$ws.Range("M34").Value = -1076
$ws.Range("N34").Value = -2992.1177

$ws.Range("H58").Value = 1108.775
$ws.Range("I58").Value = 1137.4667
$ws.Range("J58").Value = 1022.7
$ws.Range("K58").Value = 1137.4667
$ws.Range("L58").Value = 1022.7
$ws.Range("M58").Value = -934.4666999999999
$ws.Range("N58").Value = -1428.7

$ws.Range("H99").Value = 1505.7084
$ws.Range("I99").Value = 1456.5385
$ws.Range("J99").Value = 1563.8182
$ws.Range("K99").Value = 1456.5385
$ws.Range("L99").Value = 1563.8182
$ws.Range("M99").Value = 41.46149999999989
$ws.Range("N99").Value = -4559.8182

$ws.Range("H126").Value = 1505.7084
$ws.Range("I126").Value = 1456.5385
$ws.Range("J126").Value = 1563.8182
$ws.Range("K126").Value = 4369.6155
$ws.Range("L126").Value = 4691.4546
$ws.Range("M126").Value = -1899.6155
$ws.Range("N126").Value = -9631.454600000001

$ws.Range("H132").Value = 2197.325
$ws.Range("I132").Value = 1643.8823
$ws.Range("J132").Value = 5333.5
$ws.Range("K132").Value = 4931.6469
$ws.Range("L132").Value = 16000.5
$ws.Range("M132").Value = -2401.6469
$ws.Range("N132").Value = -21060.5

$ws.Range("H134").Value = 4014.9565
$ws.Range("I134").Value = 4627.2905
$ws.Range("K134").Value = 13881.8715
$ws.Range("M134").Value = -11346.8715

$ws.Range("H136").Value = 1108.775
$ws.Range("I136").Value = 1137.4667
$ws.Range("J136").Value = 1022.7
$ws.Range("K136").Value = 3412.4001
$ws.Range("L136").Value = 3068.1
$ws.Range("M136").Value = -862.4000999999998
$ws.Range("N136").Value = -8168.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1874
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1874
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 5622
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -8118

$ws.Range("H94").Value = 6067.6665
$ws.Range("I94").Value = 10000
$ws.Range("J94").Value = 5786.7856
$ws.Range("K94").Value = 30000
$ws.Range("L94").Value = 17360.3568
$ws.Range("M94").Value = -29324
$ws.Range("N94").Value = -18712.3568

$ws.Range("H96").Value = 4084.6155
$ws.Range("J96").Value = 4084.6155
$ws.Range("L96").Value = 12253.8465
$ws.Range("N96").Value = -16371.8465

$ws.Range("H123").Value = 3043.3333
$ws.Range("I123").Value = 1915
$ws.Range("J123").Value = 5300
$ws.Range("K123").Value = 5745
$ws.Range("L123").Value = 15900
$ws.Range("M123").Value = -3295
$ws.Range("N123").Value = -20800

$ws.Range("H129").Value = 1384
$ws.Range("J129").Value = 1513.56
$ws.Range("L129").Value = 4540.68
$ws.Range("N129").Value = -14540.68

$ws.Range("H131").Value = 908.08
$ws.Range("I131").Value = 458.33334
$ws.Range("J131").Value = 936.78723
$ws.Range("K131").Value = 1375.00002
$ws.Range("L131").Value = 2810.36169
$ws.Range("M131").Value = 3664.99998
$ws.Range("N131").Value = -12890.36169

$ws.Range("H133").Value = 1860
$ws.Range("I133").Value = 1860
$ws.Range("K133").Value = 5580
$ws.Range("M133").Value = -520

$ws.Range("H134").Value = 2462.7334
$ws.Range("I134").Value = 1762.8182
$ws.Range("J134").Value = 4387.5
$ws.Range("K134").Value = 5288.4546
$ws.Range("L134").Value = 13162.5
$ws.Range("M134").Value = -218.4546
$ws.Range("N134").Value = -23302.5

$ws.Range("H136").Value = 4174.0146
$ws.Range("I136").Value = 2190.8333
$ws.Range("J136").Value = 4591.5264
$ws.Range("K136").Value = 6572.499899999999
$ws.Range("L136").Value = 13774.5792
$ws.Range("M136").Value = -1472.499899999999
$ws.Range("N136").Value = -23974.5792

$ws.Range("H137").Value = 2006.6666
$ws.Range("I137").Value = 1831.5385
$ws.Range("J137").Value = 2291.25
$ws.Range("K137").Value = 5494.6155
$ws.Range("L137").Value = 6873.75
$ws.Range("M137").Value = -394.6154999999999
$ws.Range("N137").Value = -17073.75

$ws.Range("H138").Value = 29414060
$ws.Range("I138").Value = 76924210
$ws.Range("J138").Value = 3014.2856
$ws.Range("K138").Value = 230772630
$ws.Range("L138").Value = 9042.856800000001
$ws.Range("M138").Value = -230767490
$ws.Range("N138").Value = -19322.8568

$ws.Range("H139").Value = 23866.695
$ws.Range("I139").Value = 1568.1515
$ws.Range("J139").Value = 80470.69500000001
$ws.Range("K139").Value = 4704.4545
$ws.Range("L139").Value = 241412.085
$ws.Range("M139").Value = 435.5455000000002
$ws.Range("N139").Value = -251692.085

$ws.Range("H140").Value = 33552.605
$ws.Range("I140").Value = 63765.5
$ws.Range("K140").Value = 191296.5
$ws.Range("M140").Value = -186116.5

$ws.Range("H141").Value = 10667.4375
$ws.Range("I141").Value = 4408.778
$ws.Range("J141").Value = 18714.285
$ws.Range("K141").Value = 13226.334
$ws.Range("L141").Value = 56142.855
$ws.Range("M141").Value = -8046.334000000001
$ws.Range("N141").Value = -66502.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1934.8889
$ws.Range("I132").Value = 1455.5405
$ws.Range("J132").Value = 4151.875
$ws.Range("K132").Value = 4366.6215
$ws.Range("L132").Value = 12455.625
$ws.Range("M132").Value = -1836.6215
$ws.Range("N132").Value = -17515.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1969.3928
$ws.Range("I132").Value = 1799.3721
$ws.Range("J132").Value = 2531.7693
$ws.Range("K132").Value = 5398.1163
$ws.Range("L132").Value = 7595.3079
$ws.Range("M132").Value = -2868.1163
$ws.Range("N132").Value = -12655.3079

$ws.Range("H136").Value = 1396.6792
$ws.Range("I136").Value = 1296.2
$ws.Range("J136").Value = 1961.875
$ws.Range("K136").Value = 3888.6
$ws.Range("L136").Value = 5885.625
$ws.Range("M136").Value = -1338.6
$ws.Range("N136").Value = -10985.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1693.2898
$ws.Range("I132").Value = 2013.7778
$ws.Range("J132").Value = 1092.375
$ws.Range("K132").Value = 6041.3334
$ws.Range("L132").Value = 3277.125
$ws.Range("M132").Value = -3511.3334
$ws.Range("N132").Value = -8337.125

$ws.Range("H136").Value = 1639.3773
$ws.Range("I136").Value = 1438.1277
$ws.Range("J136").Value = 3215.8333
$ws.Range("K136").Value = 3215.8333
$ws.Range("L136").Value = 9647.499899999999
$ws.Range("M136").Value = -1764.3831
$ws.Range("N136").Value = -14747.4999
